$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 17607
$ws.Range("B2").Value = "Sara Carvalho"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Viagem de negocios"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45096
$ws.Range("G2").Value = 3878.63

$ws.Range("A3").Value = 31270
$ws.Range("B3").Value = "João Vitor Gomes"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45099
$ws.Range("G3").Value = 3855.63

$ws.Range("A4").Value = 94321
$ws.Range("B4").Value = "Ana Carolina Araújo"
$ws.Range("C4").Value = "TI"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45106
$ws.Range("G4").Value = 6856.7

$ws.Range("A5").Value = 2144
$ws.Range("B5").Value = "Alexia Macedo"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45083
$ws.Range("G5").Value = 7253.71

$ws.Range("A6").Value = 18213
$ws.Range("B6").Value = "Maria Fernanda Carvalho"
$ws.Range("C6").Value = "Operacoes"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45096
$ws.Range("G6").Value = 5044.03

$ws.Range("A7").Value = 44621
$ws.Range("B7").Value = "Maria Cecília Sampaio"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45079
$ws.Range("G7").Value = 9759.309999999999

$ws.Range("A8").Value = 33976
$ws.Range("B8").Value = "Sra. Isabella Ferreira"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45083
$ws.Range("G8").Value = 3435.06

$ws.Range("A9").Value = 10153
$ws.Range("B9").Value = "Mirella Ferreira"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Viagem de negocios"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45093
$ws.Range("G9").Value = 5078.07

$ws.Range("A10").Value = 69028
$ws.Range("B10").Value = "João Vitor da Rosa"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45100
$ws.Range("G10").Value = 8366.969999999999

$ws.Range("A11").Value = 68374
$ws.Range("B11").Value = "Murilo Moreira"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45090
$ws.Range("G11").Value = 5688.63
